# Extend the "Demography" row (row 3) on Sheet2 with the full set of
# demography-type labels used by the simulation parameter lookup
# ("write function to evaluate abc same as cnn"): constant / shrinking /
# growing / cycling / chaotic, mirroring the existing C3:E3 "constant"
# placeholders out through columns F:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Cells.Item(3, 6).Value  = "constant"   # F3
$ws.Cells.Item(3, 7).Value  = "shrinking"  # G3
$ws.Cells.Item(3, 8).Value  = "growing"    # H3
$ws.Cells.Item(3, 9).Value  = "cycling"    # I3
$ws.Cells.Item(3, 10).Value = "chaotic"    # J3

# Move the active selection to K3, just past the newly filled range,
# matching the saved cursor position in the workbook.
$ws.Activate() | Out-Null
$ws.Range("K3").Select() | Out-Null
